# Implements "implement first version of lot sizing rules":
#  - Generic: NrBuckets 4 -> 5
#  - Productdata: updated demand-during-leadtime figures (column C/E) for the new bucket count
#  - ForecastedAverageDemand / ForcastedStandardDeviation: recomputed rows 2-5 + new row 6
#  - Capacity: recomputed aggregate capacity figures

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generic sheet: NrBuckets 4 -> 5
# ---------------------------------------------------------------------------
$wsGeneric = $wb.Worksheets.Item("Generic")
$wsGeneric.Cells.Item(4, 2).Value = 5

# ---------------------------------------------------------------------------
# Productdata sheet: columns C (r9-19) and E (r2-23) updated
# ---------------------------------------------------------------------------
$wsProd = $wb.Worksheets.Item("Productdata")

$prodE = @{
    2  = 11.7303296
    3  = 2.1915488
    4  = 0.907208
    5  = 1.452816
    6  = 0.9623999999999999
    7  = 0.2868672
    8  = 0.09757440000000001
    9  = 0.8622864
    10 = 0.5163264000000001
    11 = 0.762216
    12 = 1.2523104
    13 = 12.61465599999999
    14 = 4.889953600000001
    15 = 0.8855792
    16 = 0.8937792
    17 = 1.404
    18 = 0.449568
    19 = 0.1350272
    20 = 63.36406400000001
    21 = 67.2670976
    22 = 83.2647168
    23 = 256.2028352
}
foreach ($row in $prodE.Keys) {
    $wsProd.Cells.Item($row, 5).Value = $prodE[$row]
}

$prodC = @{
    9  = 927
    10 = 658
    11 = 2018
    12 = 477
    13 = 3742
    14 = 1396
    15 = 345
    16 = 493
    17 = 697
    18 = 198
    19 = 71
}
foreach ($row in $prodC.Keys) {
    $wsProd.Cells.Item($row, 3).Value = $prodC[$row]
}

# ---------------------------------------------------------------------------
# ForecastedAverageDemand sheet: rows 2-5 revised, row 6 appended
# columns: A..W = 1..23 ; B,E,M..W stay 0 throughout
# ---------------------------------------------------------------------------
$wsAvg = $wb.Worksheets.Item("ForecastedAverageDemand")

# row -> @{col letter index -> value} for columns C,D,F,G,H,I,J,K,L (B/E unchanged = 0)
$avgRows = @{
    2 = @{3=602; 4=157; 6=296; 7=91;  8=38; 9=419; 10=298; 11=919; 12=224}
    3 = @{3=597; 4=146; 6=297; 7=97;  8=20; 9=420; 10=298; 11=907; 12=208}
    4 = @{3=606; 4=148; 6=305; 7=91;  8=30; 9=422; 10=302; 11=901; 12=215}
    5 = @{3=605; 4=152; 6=300; 7=81;  8=32; 9=418; 10=295; 11=908; 12=213}
}
foreach ($row in $avgRows.Keys) {
    $cols = $avgRows[$row]
    foreach ($col in $cols.Keys) {
        $wsAvg.Cells.Item($row, $col).Value = $cols[$col]
    }
}

# New row 6 - clone formatting of row 5's A cell (bold/border/centered style), then fill values
$wsAvg.Range("A5").Copy() | Out-Null
$wsAvg.Range("A6").PasteSpecial(-4122) | Out-Null

$avgRow6 = @(4, 0, 607, 148, 0, 302, 86, 34, 418, 301, 902, 202, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $avgRow6.Length; $i++) {
    $wsAvg.Cells.Item(6, $i + 1).Value = $avgRow6[$i]
}

# ---------------------------------------------------------------------------
# ForcastedStandardDeviation sheet: rows 2-5 revised, row 6 appended
# ---------------------------------------------------------------------------
$wsStd = $wb.Worksheets.Item("ForcastedStandardDeviation")

$stdRows = @{
    2 = @{3=75.25;     4=19.625;   6=37;      7=4.75;      9=52.375;  10=37.25;     11=114.875;  12=28}
    3 = @{3=111.9375;  4=27.375;   6=55.6875; 7=18.1875; 8=3.75; 9=78.75;   10=55.875;    11=170.0625; 12=39}
    4 = @{3=132.5625;  7=19.90625; 8=6.5625;  9=92.3125;  10=66.0625;   11=197.09375; 12=47.03125}
    5 = @{3=141.796875; 4=35.625;  6=70.3125; 7=18.984375; 9=97.96875; 10=69.140625; 11=212.8125; 12=49.921875}
}
foreach ($row in $stdRows.Keys) {
    $cols = $stdRows[$row]
    foreach ($col in $cols.Keys) {
        $wsStd.Cells.Item($row, $col).Value = $cols[$col]
    }
}

$wsStd.Range("A5").Copy() | Out-Null
$wsStd.Range("A6").PasteSpecial(-4122) | Out-Null

$stdRow6 = @(4, 0, 147.0078125, 35.84375, 0, 73.140625, 20.828125, 8.234375, 101.234375, 72.8984375, 218.453125, 48.921875, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $stdRow6.Length; $i++) {
    $wsStd.Cells.Item(6, $i + 1).Value = $stdRow6[$i]
}

# ---------------------------------------------------------------------------
# Capacity sheet: recomputed aggregate capacity figures
# ---------------------------------------------------------------------------
$wsCap = $wb.Worksheets.Item("Capacity")
$wsCap.Cells.Item(2, 2).Value = 50274
$wsCap.Cells.Item(3, 2).Value = 527030
$wsCap.Cells.Item(4, 2).Value = 527030
$wsCap.Cells.Item(5, 2).Value = 4969140
